$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Append root-cause follow-up note to the "cause of bug" cell (C20) ---
$suffix = "不对，根本原因不在这，应该还是指令堆积在了缓冲里没来得及取出来。。。。"
$c20 = $ws.Range("C20")
$c20.Value = $c20.Value() + $suffix

# --- The "status" cell for that row no longer applies; clear it (D20) ---
$ws.Range("D20").ClearContents()

# --- Row 20 grew taller to fit the longer note ---
$ws.Rows.Item(20).RowHeight = 123

# --- Update the saved scroll position / selection for the sheet ---
$ws.Range("B20").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
